$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57:110 down to 58:111
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new weekly record
$ws.Range("A57").Value = 4
$ws.Range("B57").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C57").Value = "Los Lagos"
$ws.Range("D57").Value = 44880
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value = 10
$ws.Range("F57").Value = 100112031
$ws.Range("G57").Value = "Poroto verde"
$ws.Range("H57").Value = "Magnum"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 30
$ws.Range("K57").Value = 39000
$ws.Range("L57").Value = 39000
$ws.Range("M57").Value = 39000
$ws.Range("N57").Value = "`$/malla 25 kilos"
$ws.Range("O57").Value = "Perú"
$ws.Range("P57").Value = 1560
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"
